$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.710887666666667
$ws.Range("H2").Value = 11.132663
$ws.Range("M2").Value = 0.789222
$ws.Range("N2").Value = 2.367666
$ws.Range("O2").Value = 0.01341929863527565
$ws.Range("P2").Value = 0.01341929863527565
$ws.Range("Q2").Value = 2.928714186062
$ws.Range("R2").Value = 26.358427674558
$ws.Range("S2").Value = 0.01341929863527565
$ws.Range("T2").Value = 0.01341929863527565

# Row 3
$ws.Range("G3").Value = 3.710887666666667
$ws.Range("H3").Value = 11.132663
$ws.Range("O3").Value = 0.005047365584441773
$ws.Range("P3").Value = 0.005047365584441773
$ws.Range("Q3").Value = 1.101569582074667
$ws.Range("R3").Value = 9.914126238672001
$ws.Range("S3").Value = 0.005047365584441773
$ws.Range("T3").Value = 0.005047365584441773

# Row 4
$ws.Range("G4").Value = 3.710887666666667
$ws.Range("H4").Value = 11.132663
$ws.Range("M4").Value = 57.61405833333333
$ws.Range("N4").Value = 172.842175
$ws.Range("O4").Value = 0.9796232927683105
$ws.Range("P4").Value = 0.9796232927683105
$ws.Range("Q4").Value = 213.7992984957806
$ws.Range("R4").Value = 1924.193686462025
$ws.Range("S4").Value = 0.9796232927683105
$ws.Range("T4").Value = 0.9796232927683105

# Row 5
$ws.Range("G5").Value = 3.710887666666667
$ws.Range("H5").Value = 11.132663
$ws.Range("M5").Value = 0.1123343333333333
$ws.Range("N5").Value = 0.337003
$ws.Range("O5").Value = 0.001910043011972043
$ws.Range("P5").Value = 0.001910043011972043
$ws.Range("Q5").Value = 0.4168600921098889
$ws.Range("R5").Value = 3.751740828989
$ws.Range("S5").Value = 0.001910043011972043
$ws.Range("T5").Value = 0.001910043011972043
